$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain-text values in this workbook (e.g.
# "46.880.60"). Some of the updated prices look like ordinary decimal numbers
# (e.g. "301.42"), so Excel would otherwise auto-convert them to numeric cells
# when assigned through .Value. Force those specific cells to Text format first
# so they stay text, matching the rest of the column.
$textFormatRows = 5,6,7,9,10,11,12,16,18,19,21,22,23,25,27,28,29,30,31,32,33,34,36,38,40,41,42,43,45,48,49,50,51
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "46.843.50"
$ws.Range("E2").Value = "  +4.13%  "
$ws.Range("D3").Value = "2.259.83"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "301.42"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "99.72"
$ws.Range("E6").Value = "  +5.54%  "
$ws.Range("D7").Value = "0.560"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "35.59"
$ws.Range("E10").Value = "  +4.03%  "
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").Value = "7.16"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D14").Value = "2.605.93"
$ws.Range("D15").Value = "2.261.28"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "13.56"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "46.783.92"
$ws.Range("E17").Value = "  +4.20%  "
$ws.Range("D18").Value = "0.793"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "5.83"
$ws.Range("E21").Value = "  -3.56%  "
$ws.Range("D22").Value = "65.10"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "248.05"
$ws.Range("E23").Value = "  +3.80%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").Value = "42.14"
$ws.Range("E27").Value = "  +2.50%  "
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Value = "19.85"
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "2.79"
$ws.Range("E31").Value = "  +9.39%  "
$ws.Range("D32").Value = "145.52"
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("D33").Value = "5.40"
$ws.Range("E33").Value = "  -2.86%  "
$ws.Range("D34").Value = "0.0772"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("E35").Value = "  +9.04%  "
$ws.Range("D36").Value = "0.115"
$ws.Range("E36").Value = "  +10.40%  "
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").Value = "16.25"
$ws.Range("E38").Value = "  +19.56%  "
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("D40").Value = "3.81"
$ws.Range("E40").Value = "  -4.22%  "
$ws.Range("D41").Value = "0.0297"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").Value = "3.19"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("D45").Value = "91.55"
$ws.Range("E45").Value = "  +19.30%  "
$ws.Range("D46").Value = "1.766.70"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "71.10"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "7.84"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "93.69"
$ws.Range("E51").Value = "  -2.15%  "

# Clear the temporary Text number format again so the cells keep the workbook
# default (unstyled) appearance once the text values have been written.
foreach ($r in $textFormatRows) {
    $ws.Cells.Item($r, 4).Style = "Normal"
}
